# Deckblatt logo position center
#
# The cover-sheet ("Deckblatt") logo is an anchored picture in the only
# paragraph of the document. The edit:
#   - marks the run as "do not spell check" (w:noProof), which Word adds
#     whenever a drawing is inserted/touched,
#   - re-centers the picture horizontally (wp:posOffset -> wp:align=center),
#   - snaps the picture's stored size from 2190721x2072639 EMU to the
#     "nice" 2188800x2073600 EMU (6.08cm x 5.76cm) that Word's Format
#     Picture > Size dialog produces, updating both wp:extent and the
#     inner a:xfrm/a:ext,
#   - updates the effectExtent bounding box and the anchor's editId that
#     Word regenerates whenever an anchored object is repositioned/resized.
#
# None of the last few (editId / effectExtent) are reachable through the
# normal Shape.* property surface, so the whole paragraph is rewritten
# via Range.InsertXML with the exact target markup (this also picks up
# noProof/align/extent/ext "for free" since it is a single atomic
# replace of the run + drawing).

$d = $word.ActiveDocument

$targetParagraphXml = @'
<w:p w14:paraId="4712C575" w14:textId="7C85DD81" w:rsidR="00A740EC" w:rsidRPr="00EB38E2" w:rsidRDefault="00FD3BA2" w:rsidP="00E3648D" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:spacing w:line="600" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r w:rsidRPr="00FD3BA2"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="368FA46F" wp14:editId="519A1046"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>3062605</wp:posOffset></wp:positionV><wp:extent cx="2188800" cy="2073600"/><wp:effectExtent l="0" t="0" r="2540" b="3175"/><wp:wrapTopAndBottom/><wp:docPr id="1303869125" name="Grafik 1" descr="Ein Bild, das Emblem, Symbol, Markenzeichen, Kreis enth&#228;lt.&#10;&#10;Automatisch generierte Beschreibung"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="1303869125" name="Grafik 1" descr="Ein Bild, das Emblem, Symbol, Markenzeichen, Kreis enth&#228;lt.&#10;&#10;Automatisch generierte Beschreibung"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2188800" cy="2073600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>
'@

$targetRange = $d.Paragraphs(1).Range
$targetRange.InsertXML($targetParagraphXml)
